$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new login-log row (row 24) mirroring the existing entries in
# columns A-D: data_login, hora_login, usuario, email.
#
# Column A holds dates stored as plain text (shared strings) in this sheet,
# so force a text number format before assigning the value to stop Excel
# from auto-converting "2026-01-06" into a date serial, then drop the
# number format back to Normal/General so the cell keeps the same
# (unstyled) look as the other data rows.
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "2026-01-06"
$ws.Range("A24").Style = "Normal"

$ws.Range("B24").Value = "17:09:45"
$ws.Range("C24").Value = "maria"
$ws.Range("D24").Value = "maria@teste.com"
